$d = $word.ActiveDocument

$d.Content.Find.Execute("Concluído (91 dias)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Concluído (95 dias)", 2)

$d.Content.Find.Execute("Assinatura Contrato (100 dias)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Assinatura Contrato (104 dias)", 2)

$d.Content.Find.Execute("Assinatura Contrato (77 dias)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Assinatura Contrato (81 dias)", 2)

$d.Content.Find.Execute("Assinatura Contrato (63 dias)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Assinatura Contrato (67 dias)", 2)

$d.Content.Find.Execute("AGU (11 dias)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "AGU (15 dias)", 2)

$d.Content.Find.Execute("Total de dias 678", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Total de dias 698", 2)
